$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B..F to C..G)
$ws.Columns("B:B").Insert()

# Move the segment-name strings from column A to the new column B,
# and put the numeric segment index (0-based) into column A.
$segmentNames = @("background","back_bumper","back_glass","back_left_door","back_left_light","back_right_door","back_right_light","front_bumper","front_glass","front_left_door","front_left_light","front_right_door","front_right_light","hood","left_mirror","right_mirror","tailgate","trunk","wheel")

for ($i = 0; $i -lt $segmentNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $segmentNames[$i]
    $ws.Cells.Item($row, 1).Value = $i
}

# New column B data cells should have no special style (same as other data columns)
$ws.Range("B2:B20").ClearFormats()

# Set the new header for column B (same bold/centered style as the other headers)
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"
